$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Reorder slides.
#    The "mhi" scatter-plot slide (currently slide #11) moves to the very
#    end of the deck, and a brand-new "Questions and Discussion" section
#    header slide is inserted right before the final ("mhi formulas") slide.
# ---------------------------------------------------------------------------
$p.Slides.Item(11).MoveTo($p.Slides.Count)

$newSlide = $p.Slides.Add($p.Slides.Count, 12)
$newSlide.CustomLayout = $p.SlideMaster.CustomLayouts.Item(3)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Questions and Discussion"

# ---------------------------------------------------------------------------
# 2. Edit the "mhi formulas" slide (the explanatory slide with the rho/sigma
#    formulas) - rename the "mhi" variable to "hhsize" and re-flow a couple
#    of boxes/connector that sit next to it.
# ---------------------------------------------------------------------------
$formulaSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.Count -gt 0 -and $candidate.Shapes.Item(1).Name -eq "TextBox 1") {
        $txt = $candidate.Shapes.Item(1).TextFrame.TextRange.Text
        if ($txt -like "*pov*") {
            $formulaSlide = $candidate
        }
    }
}

foreach ($shape in $formulaSlide.Shapes) {
    if ($shape.Name -eq "TextBox 7") {
        $shape.Left = 7139853 / 12700.0
        $shape.Width = 996981 / 12700.0
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("mhi")
        $tr.Characters($idx + 1, 3).Text = "hhsize"
    }
    elseif ($shape.Name -eq "TextBox 8") {
        $shape.Left = 8770850 / 12700.0
    }
    elseif ($shape.Name -eq "Straight Arrow Connector 9") {
        $shape.Left = 8136834 / 12700.0
        $shape.Width = 634016 / 12700.0
    }
    elseif ($shape.Name -eq "TextBox 11") {
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("mhi")
        $tr.Characters($idx + 1, 3).Text = "hhsize"
        $full2 = $tr.Text
        $idx2 = $full2.IndexOf(" = median household income")
        $tr.Characters($idx2 + 1, " = median household income".Length).Text = " = household size"
    }
}

# ---------------------------------------------------------------------------
# 3. Edit the "Why Springfield, Missouri?" slide - reflow the bullet textbox,
#    tweak the bullet wording, and shift the map picture.
# ---------------------------------------------------------------------------
$springfieldSlide = $p.Slides.Item(2)
foreach ($shape in $springfieldSlide.Shapes) {
    if ($shape.Name -eq "TextBox 1") {
        $shape.Left = 4863548 / 12700.0
        $shape.Width = 7273305 / 12700.0
        $tr = $shape.TextFrame.TextRange
        $target = "Slave state that remained in the Union"
        $full = $tr.Text
        $idx = $full.IndexOf($target)
        $tr.Characters($idx + 1, $target.Length).Text = "Within a slave state that remained in the Union"
    }
    elseif ($shape.Name -eq "Picture 3") {
        $shape.Left = 426053 / 12700.0
    }
}
